$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.870.29"

$ws.Range("D3").Value = "1.866.66"
$ws.Range("E3").Value = "  -0.14%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5084"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3655"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.36%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07175"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8904"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.55"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.03%  "

$ws.Range("D12").Value = "1.877.28"
$ws.Range("E12").Value = "  +0.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07494"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.36%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.71%  "

$ws.Range("E16").Value = "  +0.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008487"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.17%  "

$ws.Range("E19").Value = "  +0.06%  "

$ws.Range("D20").Value = "26.915.77"
$ws.Range("E20").Value = "  -1.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.13%  "

$ws.Range("D22").Value = "2.118.68"
$ws.Range("E22").Value = "  +0.21%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.363"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.70%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.63%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.777"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.60%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.094"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.55%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.57%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.691"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.61%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.696"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09107"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.63%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05028"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7500"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.12%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.948"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.40%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.151"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.207"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.84%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.506"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01986"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.29%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5532"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.10%  "

$ws.Range("E41").Value = "  -0.11%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.569"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.56%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.30%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.581"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1484"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4738"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.33%  "

$ws.Range("E47").Value = "  +0.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.87%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.552"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.79%  "

$ws.Range("E51").Value = "  -1.12%  "

